# correção nos dados e inicio da analise PNAD 2009
#
# The sheet ("agressao/cv126101a") listed a grouping-header row before each
# block of categories (sexo, cor ou raça, grupos de idade, nível de
# instrução, classes de rendimento mensal domiciliar per capita) as well as
# an empty "sem rendimento a menos" row and the trailing source/footnote
# rows — none of which carried any data. This edit removes those rows
# (shifting the data rows below them up) and fixes the second header row
# (row 2), whose first data column had picked up the pandas placeholder
# label "unnamed: 1_level_1" instead of repeating "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the empty "section header" rows, bottom-to-top so row numbers of
# rows still to be deleted are not invalidated by the preceding deletions.
$ws.Range("A36").EntireRow.Delete()   # (1) inclusive as pessoas de cor ...
$ws.Range("A35").EntireRow.Delete()   # fonte: ibge, diretoria de pesquisas...
$ws.Range("A29").EntireRow.Delete()   # sem rendimento a menos
$ws.Range("A27").EntireRow.Delete()   # classes de rendimento mensal domiciliar per capita
$ws.Range("A19").EntireRow.Delete()   # nível de instrução
$ws.Range("A13").EntireRow.Delete()   # grupos de idade
$ws.Range("A8").EntireRow.Delete()    # cor ou raça
$ws.Range("A5").EntireRow.Delete()    # sexo

# Fix the second header row: it should mirror the "total" label used by
# row 1 instead of the stray "unnamed: 1_level_1" placeholder.
$ws.Range("B2").Value = "total"
$ws.Range("C2").Value = "própria residência"
$ws.Range("D2").Value = "via pública"
$ws.Range("E2").Value = "outro"
